# Refresh the crypto price / 1h-volume table with the latest scrape.
# (GitHub Actions "Updated cryptos list" job.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Force the cell to store $text verbatim as text even when it looks
    # like a number (e.g. "579.95" or "0.140"). Excel would otherwise
    # coerce a numeric-looking string into a real number on assignment,
    # silently dropping significant trailing zeros / changing layout.
    # The leading apostrophe is the standard Excel "treat as text" marker
    # and is stripped from the stored value.
    $ws.Range($cellRef).Value = "'" + $text
}

$ws.Range('D2').Value = '61.950.88'
$ws.Range('E2').Value = '  +4.54%  '

$ws.Range('D3').Value = '3.083.22'
$ws.Range('E3').Value = '  +2.88%  '

Set-TextValue 'D5' '579.95'
$ws.Range('E5').Value = '  +3.06%  '

Set-TextValue 'D6' '142.44'
$ws.Range('E6').Value = '  +2.38%  '

$ws.Range('E7').Value = '  +0.13%  '

$ws.Range('D8').Value = '3.068.68'
$ws.Range('E8').Value = '  +2.49%  '

Set-TextValue 'D9' '0.527'
$ws.Range('E9').Value = '  +1.32%  '

Set-TextValue 'D10' '0.140'
$ws.Range('E10').Value = '  +4.69%  '

Set-TextValue 'D11' '5.58'
$ws.Range('E11').Value = '  +8.81%  '

Set-TextValue 'D12' '0.467'
$ws.Range('E12').Value = '  +2.53%  '

Set-TextValue 'D13' '0.0000240'
$ws.Range('E13').Value = '  +3.85%  '

Set-TextValue 'D14' '35.34'
$ws.Range('E14').Value = '  +4.76%  '

$ws.Range('E15').Value = '  +0.18%  '

$ws.Range('D16').Value = '3.591.62'
$ws.Range('E16').Value = '  +2.92%  '

Set-TextValue 'D17' '7.28'
$ws.Range('E17').Value = '  +2.94%  '

$ws.Range('D18').Value = '3.076.96'
$ws.Range('E18').Value = '  +2.82%  '

$ws.Range('D19').Value = '61.850.23'
$ws.Range('E19').Value = '  +4.51%  '

Set-TextValue 'D20' '449.75'
$ws.Range('E20').Value = '  +4.72%  '

Set-TextValue 'D21' '13.92'
$ws.Range('E21').Value = '  +2.28%  '

Set-TextValue 'D22' '0.729'
$ws.Range('E22').Value = '  +1.99%  '

Set-TextValue 'D23' '7.43'
$ws.Range('E23').Value = '  +4.68%  '

Set-TextValue 'D24' '13.78'
$ws.Range('E24').Value = '  +2.70%  '

Set-TextValue 'D25' '81.95'
$ws.Range('E25').Value = '  +1.48%  '

$ws.Range('E26').Value = '  +0.00%  '

Set-TextValue 'D27' '2.27'
$ws.Range('E27').Value = '  +5.10%  '

$ws.Range('E28').Value = '  +0.07%  '

Set-TextValue 'D29' '2.66'
$ws.Range('E29').Value = '  +4.78%  '

Set-TextValue 'D30' '8.14'
$ws.Range('E30').Value = '  +4.88%  '

Set-TextValue 'D31' '6.76'

Set-TextValue 'D32' '0.111'
$ws.Range('E32').Value = '  +12.39%  '

Set-TextValue 'D33' '26.76'

$ws.Range('E34').Value = '  +3.92%  '

$ws.Range('D35').Value = '0.0₃0801'
$ws.Range('E35').Value = '  +2.91%  '

Set-TextValue 'D36' '6.05'
$ws.Range('E36').Value = '  +3.60%  '

Set-TextValue 'D37' '2.20'
$ws.Range('E37').Value = '  +6.05%  '

Set-TextValue 'D38' '50.37'
$ws.Range('E38').Value = '  +2.28%  '

Set-TextValue 'D39' '3.01'
$ws.Range('E39').Value = '  +9.01%  '

Set-TextValue 'D40' '8.84'
$ws.Range('E40').Value = '  +2.23%  '

Set-TextValue 'D41' '430.53'
$ws.Range('E41').Value = '  +6.73%  '

Set-TextValue 'D42' '0.0372'
$ws.Range('E42').Value = '  +5.84%  '

$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue 'D43' '0.273'
$ws.Range('E43').Value = '  +7.81%  '

$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.797.74'
$ws.Range('E44').Value = '  +1.24%  '

$ws.Range('E45').Value = '  +0.93%  '

Set-TextValue 'D46' '2.12'
$ws.Range('E46').Value = '  +5.65%  '

$ws.Range('B47').Value = 'Arweave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextValue 'D47' '35.23'
$ws.Range('E47').Value = '  +4.80%  '

$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D48' '0.999'
$ws.Range('E48').Value = '  -0.02%  '

Set-TextValue 'D49' '123.77'
$ws.Range('E49').Value = '  +0.58%  '

$ws.Range('E50').Value = '  +1.01%  '

Set-TextValue 'D51' '24.09'
$ws.Range('E51').Value = '  +2.38%  '
